$d = $word.ActiveDocument

$d.Content.Find.Execute("85×45=", $true, $false, $false, $false, $false, $true, 1, $false, "78×44=", 2) | Out-Null
$d.Content.Find.Execute("56×21=", $true, $false, $false, $false, $false, $true, 1, $false, "83×36=", 2) | Out-Null
$d.Content.Find.Execute("48×27=", $true, $false, $false, $false, $false, $true, 1, $false, "14×10=", 2) | Out-Null
$d.Content.Find.Execute("71×37=", $true, $false, $false, $false, $false, $true, 1, $false, "99×58=", 2) | Out-Null
$d.Content.Find.Execute("27×64=", $true, $false, $false, $false, $false, $true, 1, $false, "35×65=", 2) | Out-Null
$d.Content.Find.Execute("15×72=", $true, $false, $false, $false, $false, $true, 1, $false, "40×57=", 2) | Out-Null
$d.Content.Find.Execute("97×31=", $true, $false, $false, $false, $false, $true, 1, $false, "95×99=", 2) | Out-Null
$d.Content.Find.Execute("98×39=", $true, $false, $false, $false, $false, $true, 1, $false, "43×36=", 2) | Out-Null
$d.Content.Find.Execute("31×92=", $true, $false, $false, $false, $false, $true, 1, $false, "32×62=", 2) | Out-Null
$d.Content.Find.Execute("69×49=", $true, $false, $false, $false, $false, $true, 1, $false, "22×87=", 2) | Out-Null
$d.Content.Find.Execute("40×96=", $true, $false, $false, $false, $false, $true, 1, $false, "28×64=", 2) | Out-Null
$d.Content.Find.Execute("95×62=", $true, $false, $false, $false, $false, $true, 1, $false, "65×98=", 2) | Out-Null
$d.Content.Find.Execute("63×55=", $true, $false, $false, $false, $false, $true, 1, $false, "76×97=", 2) | Out-Null
$d.Content.Find.Execute("78×55=", $true, $false, $false, $false, $false, $true, 1, $false, "91×80=", 2) | Out-Null
$d.Content.Find.Execute("67×85=", $true, $false, $false, $false, $false, $true, 1, $false, "39×92=", 2) | Out-Null
$d.Content.Find.Execute("79×70=", $true, $false, $false, $false, $false, $true, 1, $false, "46×48=", 2) | Out-Null
$d.Content.Find.Execute("46×23=", $true, $false, $false, $false, $false, $true, 1, $false, "57×32=", 2) | Out-Null
$d.Content.Find.Execute("83×64=", $true, $false, $false, $false, $false, $true, 1, $false, "99×66=", 2) | Out-Null
$d.Content.Find.Execute("12×56=", $true, $false, $false, $false, $false, $true, 1, $false, "29×15=", 2) | Out-Null
$d.Content.Find.Execute("46×10=", $true, $false, $false, $false, $false, $true, 1, $false, "44×29=", 2) | Out-Null
$d.Content.Find.Execute("74×96=", $true, $false, $false, $false, $false, $true, 1, $false, "82×61=", 2) | Out-Null
$d.Content.Find.Execute("18×42=", $true, $false, $false, $false, $false, $true, 1, $false, "36×17=", 2) | Out-Null
$d.Content.Find.Execute("20×67=", $true, $false, $false, $false, $false, $true, 1, $false, "86×91=", 2) | Out-Null
$d.Content.Find.Execute("85×97=", $true, $false, $false, $false, $false, $true, 1, $false, "94×85=", 2) | Out-Null
$d.Content.Find.Execute("45×56=", $true, $false, $false, $false, $false, $true, 1, $false, "63×31=", 2) | Out-Null
$d.Content.Find.Execute("93×87=", $true, $false, $false, $false, $false, $true, 1, $false, "95×81=", 2) | Out-Null
$d.Content.Find.Execute("52×17=", $true, $false, $false, $false, $false, $true, 1, $false, "56×66=", 2) | Out-Null
$d.Content.Find.Execute("31×39=", $true, $false, $false, $false, $false, $true, 1, $false, "69×22=", 2) | Out-Null
$d.Content.Find.Execute("92×56=", $true, $false, $false, $false, $false, $true, 1, $false, "64×67=", 2) | Out-Null
$d.Content.Find.Execute("34×53=", $true, $false, $false, $false, $false, $true, 1, $false, "94×68=", 2) | Out-Null
$d.Content.Find.Execute("58×40=", $true, $false, $false, $false, $false, $true, 1, $false, "92×66=", 2) | Out-Null
$d.Content.Find.Execute("44×39=", $true, $false, $false, $false, $false, $true, 1, $false, "34×83=", 2) | Out-Null
$d.Content.Find.Execute("96×12=", $true, $false, $false, $false, $false, $true, 1, $false, "25×28=", 2) | Out-Null
$d.Content.Find.Execute("72×12=", $true, $false, $false, $false, $false, $true, 1, $false, "75×83=", 2) | Out-Null
$d.Content.Find.Execute("71×69=", $true, $false, $false, $false, $false, $true, 1, $false, "56×46=", 2) | Out-Null
$d.Content.Find.Execute("69×80=", $true, $false, $false, $false, $false, $true, 1, $false, "86×43=", 2) | Out-Null
$d.Content.Find.Execute("24×12=", $true, $false, $false, $false, $false, $true, 1, $false, "46×68=", 2) | Out-Null
$d.Content.Find.Execute("17×71=", $true, $false, $false, $false, $false, $true, 1, $false, "86×72=", 2) | Out-Null
$d.Content.Find.Execute("70×83=", $true, $false, $false, $false, $false, $true, 1, $false, "79×38=", 2) | Out-Null
$d.Content.Find.Execute("32×41=", $true, $false, $false, $false, $false, $true, 1, $false, "95×42=", 2) | Out-Null
$d.Content.Find.Execute("73×83=", $true, $false, $false, $false, $false, $true, 1, $false, "15×79=", 2) | Out-Null
$d.Content.Find.Execute("95×54=", $true, $false, $false, $false, $false, $true, 1, $false, "30×12=", 2) | Out-Null
$d.Content.Find.Execute("36×20=", $true, $false, $false, $false, $false, $true, 1, $false, "40×30=", 2) | Out-Null
$d.Content.Find.Execute("34×65=", $true, $false, $false, $false, $false, $true, 1, $false, "90×26=", 2) | Out-Null
$d.Content.Find.Execute("43×77=", $true, $false, $false, $false, $false, $true, 1, $false, "14×92=", 2) | Out-Null
$d.Content.Find.Execute("40×58=", $true, $false, $false, $false, $false, $true, 1, $false, "51×54=", 2) | Out-Null
$d.Content.Find.Execute("92×85=", $true, $false, $false, $false, $false, $true, 1, $false, "84×26=", 2) | Out-Null
$d.Content.Find.Execute("23×60=", $true, $false, $false, $false, $false, $true, 1, $false, "53×35=", 2) | Out-Null
$d.Content.Find.Execute("35×91=", $true, $false, $false, $false, $false, $true, 1, $false, "46×21=", 2) | Out-Null
$d.Content.Find.Execute("100×75=", $true, $false, $false, $false, $false, $true, 1, $false, "29×25=", 2) | Out-Null
$d.Content.Find.Execute("73×36=", $true, $false, $false, $false, $false, $true, 1, $false, "98×92=", 2) | Out-Null
$d.Content.Find.Execute("18×97=", $true, $false, $false, $false, $false, $true, 1, $false, "83×75=", 2) | Out-Null
$d.Content.Find.Execute("57×38=", $true, $false, $false, $false, $false, $true, 1, $false, "27×24=", 2) | Out-Null
$d.Content.Find.Execute("25×91=", $true, $false, $false, $false, $false, $true, 1, $false, "77×100=", 2) | Out-Null
$d.Content.Find.Execute("26×73=", $true, $false, $false, $false, $false, $true, 1, $false, "78×92=", 2) | Out-Null
$d.Content.Find.Execute("70×79=", $true, $false, $false, $false, $false, $true, 1, $false, "75×29=", 2) | Out-Null
$d.Content.Find.Execute("59×87=", $true, $false, $false, $false, $false, $true, 1, $false, "10×36=", 2) | Out-Null
$d.Content.Find.Execute("16×77=", $true, $false, $false, $false, $false, $true, 1, $false, "10×80=", 2) | Out-Null
$d.Content.Find.Execute("35×42=", $true, $false, $false, $false, $false, $true, 1, $false, "24×76=", 2) | Out-Null
$d.Content.Find.Execute("17×23=", $true, $false, $false, $false, $false, $true, 1, $false, "53×63=", 2) | Out-Null
$d.Content.Find.Execute("27×92=", $true, $false, $false, $false, $false, $true, 1, $false, "36×46=", 2) | Out-Null
$d.Content.Find.Execute("65×38=", $true, $false, $false, $false, $false, $true, 1, $false, "65×15=", 2) | Out-Null
$d.Content.Find.Execute("98×90=", $true, $false, $false, $false, $false, $true, 1, $false, "64×88=", 2) | Out-Null
$d.Content.Find.Execute("48×35=", $true, $false, $false, $false, $false, $true, 1, $false, "43×12=", 2) | Out-Null
$d.Content.Find.Execute("72×13=", $true, $false, $false, $false, $false, $true, 1, $false, "81×67=", 2) | Out-Null
$d.Content.Find.Execute("23×75=", $true, $false, $false, $false, $false, $true, 1, $false, "68×78=", 2) | Out-Null
$d.Content.Find.Execute("28×23=", $true, $false, $false, $false, $false, $true, 1, $false, "39×10=", 2) | Out-Null
$d.Content.Find.Execute("21×78=", $true, $false, $false, $false, $false, $true, 1, $false, "61×18=", 2) | Out-Null
$d.Content.Find.Execute("58×76=", $true, $false, $false, $false, $false, $true, 1, $false, "80×76=", 2) | Out-Null
$d.Content.Find.Execute("77×99=", $true, $false, $false, $false, $false, $true, 1, $false, "86×25=", 2) | Out-Null
$d.Content.Find.Execute("33×59=", $true, $false, $false, $false, $false, $true, 1, $false, "43×91=", 2) | Out-Null
$d.Content.Find.Execute("48×98=", $true, $false, $false, $false, $false, $true, 1, $false, "97×55=", 2) | Out-Null
$d.Content.Find.Execute("44×46=", $true, $false, $false, $false, $false, $true, 1, $false, "73×20=", 2) | Out-Null
$d.Content.Find.Execute("19×39=", $true, $false, $false, $false, $false, $true, 1, $false, "70×94=", 2) | Out-Null
$d.Content.Find.Execute("56×82=", $true, $false, $false, $false, $false, $true, 1, $false, "25×84=", 2) | Out-Null
$d.Content.Find.Execute("11×98=", $true, $false, $false, $false, $false, $true, 1, $false, "44×97=", 2) | Out-Null
$d.Content.Find.Execute("25×98=", $true, $false, $false, $false, $false, $true, 1, $false, "16×21=", 2) | Out-Null
$d.Content.Find.Execute("38×37=", $true, $false, $false, $false, $false, $true, 1, $false, "43×59=", 2) | Out-Null
$d.Content.Find.Execute("13×31=", $true, $false, $false, $false, $false, $true, 1, $false, "91×57=", 2) | Out-Null
$d.Content.Find.Execute("93×93=", $true, $false, $false, $false, $false, $true, 1, $false, "30×14=", 2) | Out-Null
$d.Content.Find.Execute("78×12=", $true, $false, $false, $false, $false, $true, 1, $false, "92×43=", 2) | Out-Null
$d.Content.Find.Execute("44×91=", $true, $false, $false, $false, $false, $true, 1, $false, "36×88=", 2) | Out-Null
$d.Content.Find.Execute("89×80=", $true, $false, $false, $false, $false, $true, 1, $false, "37×88=", 2) | Out-Null
$d.Content.Find.Execute("66×16=", $true, $false, $false, $false, $false, $true, 1, $false, "41×33=", 2) | Out-Null
$d.Content.Find.Execute("22×50=", $true, $false, $false, $false, $false, $true, 1, $false, "81×57=", 2) | Out-Null
$d.Content.Find.Execute("92×95=", $true, $false, $false, $false, $false, $true, 1, $false, "50×92=", 2) | Out-Null
$d.Content.Find.Execute("59×48=", $true, $false, $false, $false, $false, $true, 1, $false, "98×56=", 2) | Out-Null
$d.Content.Find.Execute("53×65=", $true, $false, $false, $false, $false, $true, 1, $false, "88×89=", 2) | Out-Null
$d.Content.Find.Execute("92×32=", $true, $false, $false, $false, $false, $true, 1, $false, "36×41=", 2) | Out-Null
$d.Content.Find.Execute("33×80=", $true, $false, $false, $false, $false, $true, 1, $false, "63×70=", 2) | Out-Null
$d.Content.Find.Execute("91×22=", $true, $false, $false, $false, $false, $true, 1, $false, "96×54=", 2) | Out-Null
$d.Content.Find.Execute("68×61=", $true, $false, $false, $false, $false, $true, 1, $false, "66×29=", 2) | Out-Null
$d.Content.Find.Execute("59×34=", $true, $false, $false, $false, $false, $true, 1, $false, "55×68=", 2) | Out-Null
$d.Content.Find.Execute("54×66=", $true, $false, $false, $false, $false, $true, 1, $false, "67×56=", 2) | Out-Null
$d.Content.Find.Execute("46×62=", $true, $false, $false, $false, $false, $true, 1, $false, "50×66=", 2) | Out-Null
$d.Content.Find.Execute("100×58=", $true, $false, $false, $false, $false, $true, 1, $false, "55×73=", 2) | Out-Null
$d.Content.Find.Execute("45×15=", $true, $false, $false, $false, $false, $true, 1, $false, "53×92=", 2) | Out-Null
$d.Content.Find.Execute("22×75=", $true, $false, $false, $false, $false, $true, 1, $false, "95×56=", 2) | Out-Null
$d.Content.Find.Execute("75×18=", $true, $false, $false, $false, $false, $true, 1, $false, "48×22=", 2) | Out-Null
$d.Content.Find.Execute("43×13=", $true, $false, $false, $false, $false, $true, 1, $false, "79×80=", 2) | Out-Null
